# Apply the "540 collection in blockes" edit:
# - Clear a set of F-column "(531) ..." classification strings
# - Clear a set of E-column repeated header-continuation strings
# - Extend two E-column marca names with additional words

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear column F cells that held Vienna classification codes
$fClearRows = @(4, 6, 8, 9, 12, 13, 15, 17, 20, 21, 24, 34, 38, 43, 44, 46, 47, 48, 53, 56, 57, 58, 59, 61, 62, 63, 69)
foreach ($r in $fClearRows) {
    $ws.Range("F$r").Value = ""
}

# Clear column E cells that held stray repeated text
$eClearRows = @(16, 30, 39, 41, 49, 68)
foreach ($r in $eClearRows) {
    $ws.Range("E$r").Value = ""
}

# Update/extend two marca (brand) names in column E
$ws.Range("E29").Value = "ESCOLAS DE CONDUÇÃO GRUPO LIDADOR DA MAIA"
$ws.Range("E64").Value = "SAVOR THE MIND'S REBELLION"
